$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" between "2021-Q3" and "总计"
#
# NOTE: worksheet handles in this host are positional, so they must
# be re-fetched by name immediately after any structural operation
# (Add / Move / Delete) shifts tab positions around.
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Item("总计")
$wsNew = $wb.Worksheets.Item("2022-Q1")
$wsNew.Move($wsTotal)

# Re-fetch again after Move() reshuffled tab positions.
$ws = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# Reuse the existing header / index-column formatting from the "总计"
# sheet (same bold/centered/bordered style used across the workbook)
# so no extra styles are introduced.
$wsTotal.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data row — fund codes / ratios are kept as text (leading apostrophe)
# so values like "005126" and "0.10" keep their exact formatting
# instead of being coerced to numbers.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'005126"
$ws.Range("C2").Value = "银河量化稳进混合"
$ws.Range("D2").Value = "'0.10"
$ws.Range("E2").Value = "'78.20"
$ws.Range("F2").Value = "'2.56"
$ws.Range("G2").Value = "'0.0026"
$ws.Range("H2").Value = 1

# ------------------------------------------------------------------
# 2) Update the "总计" sheet: add a new first data row for 2022-Q1,
#    pushing the existing 2021-Q3 row down and renumbering the index.
# ------------------------------------------------------------------
# Move the current row 2 (2021-Q3) data down to row 3 first.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.08

# Copy the index-column format onto the newly used A3 cell.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Now write the new 2022-Q1 summary row into row 2.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0

Write-Host "Edit complete"
